$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 65 ("Perejil" weekly update),
# pushing the former rows 65-67 down to 67-69.
$ws.Rows("65:66").Insert()

# New row 65: week of 2023-05-29 (serial 45075), quality "Primera"
$ws.Cells.Item(65, 1).Value = 7
$ws.Cells.Item(65, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(65, 3).Value = "Ñuble"
$ws.Cells.Item(65, 4).Value = 45075
$ws.Cells.Item(65, 5).Value = 16
$ws.Cells.Item(65, 6).Value = 100112044
$ws.Cells.Item(65, 7).Value = "Perejil"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 150
$ws.Cells.Item(65, 11).Value = 1200
$ws.Cells.Item(65, 12).Value = 1200
$ws.Cells.Item(65, 13).Value = 1200
$ws.Cells.Item(65, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(65, 15).Value = "Región del Maule"
$ws.Cells.Item(65, 16).Value = 1200
$ws.Cells.Item(65, 17).Value = 1
$ws.Cells.Item(65, 18).Value = "Hortaliza"

# New row 66: week of 2023-05-29 (serial 45075), quality "Segunda"
$ws.Cells.Item(66, 1).Value = 7
$ws.Cells.Item(66, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(66, 3).Value = "Ñuble"
$ws.Cells.Item(66, 4).Value = 45075
$ws.Cells.Item(66, 5).Value = 16
$ws.Cells.Item(66, 6).Value = 100112044
$ws.Cells.Item(66, 7).Value = "Perejil"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Segunda"
$ws.Cells.Item(66, 10).Value = 100
$ws.Cells.Item(66, 11).Value = 1000
$ws.Cells.Item(66, 12).Value = 1000
$ws.Cells.Item(66, 13).Value = 1000
$ws.Cells.Item(66, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(66, 15).Value = "Región del Maule"
$ws.Cells.Item(66, 16).Value = 1000
$ws.Cells.Item(66, 17).Value = 1
$ws.Cells.Item(66, 18).Value = "Hortaliza"
